$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price list values (rows 2-11, columns B-H)
$data = @{
    2  = @(89003, 120334, 150444, 255048, 440277, 77516, 66864)
    3  = @(89391, 129370, 151869, 262068, 441641, 77719, 67127)
    4  = @(100556, 129635, 171127, 294974, 497270, 87340, 75717)
    5  = @(141938, 202945, 237711, 398345, 593139, 123818, 107697)
    6  = @(244056, 325094, 427876, 596382, 812707, 213848, 186636)
    7  = @(63825, 91990, 114222, 215890, 314976, 60672, 54604)
    8  = @(65023, 92898, 114185, 217669, 320823, 60382, 56357)
    9  = @(74629, 105159, 131835, 252550, 370441, 70594, 66126)
    10 = @(90953, 127857, 157028, 304054, 433478, 84809, 79441)
    11 = @(75652, 102282, 127876, 216790, 374235, 65890, 56834)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        # Column B is index 2, C is 3, ... H is 8
        $col = $i + 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

# Update the selected cell in the sheet view
$ws.Range("E18").Select()
